$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the mis-typed test names: "...Summer-2015-2018" should read "...Summer-2015-2016"
foreach ($row in 12..14) {
    $ws.Range("A$row").Value = "(SBAC_PT)SBAC-IRP-Perf-MATH-3-Summer-2015-2016"
}

foreach ($row in 36..38) {
    $ws.Range("A$row").Value = "(SBAC_PT)SBAC-IRP-Perf-MATH-11-Summer-2015-2016"
}
